# Praktikum 1 - Membuat RESTful API Register
# The template_user sheet holds a single sample "register" row:
#   level_id | username | nama | password
#      13    | Direktur | Om Direktur | 12345
# Update the sample level_id from 10 to 13 and leave the cursor on A2,
# matching the row used to exercise the new Register endpoint.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged text, re-asserted defensively)
$ws.Range("A1").Value = "level_id"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "nama"
$ws.Range("D1").Value = "password"

# Sample data row
$ws.Range("A2").Value = 13
$ws.Range("B2").Value = "Direktur"
$ws.Range("C2").Value = "Om Direktur"
$ws.Range("D2").Value = 12345

# Move the active selection to A2
$ws.Range("A2").Select()
